$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Recorded By") lists the users who recorded each attendance
# session as a comma-separated string. Whenever "System" is the first
# entry in such a list, swap it with the last entry in the list so that
# "System" is no longer reported as the primary recorder.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val) {
        $text = $val.ToString()
        if ($text.StartsWith("System, ")) {
            $parts = $text.Split(", ")
            if ($parts.Length -gt 1) {
                $first = $parts[0]
                $last = $parts[$parts.Length - 1]
                $parts[0] = $last
                $parts[$parts.Length - 1] = $first
                $newVal = [string]::Join(", ", $parts)
                $cell.Value2 = $newVal
            }
        }
    }
}
